$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 814
$ws.Range("I18").Value = 872.0909
$ws.Range("K18").Value = 872.0909
$ws.Range("M18").Value = -588.0909

$ws.Range("H33").Value = 282.6154
$ws.Range("J33").Value = 242.33333
$ws.Range("L33").Value = 242.33333
$ws.Range("N33").Value = -700.3333299999999

$ws.Range("H95").Value = 48798
$ws.Range("J95").Value = 48798
$ws.Range("L95").Value = 48798
$ws.Range("N95").Value = -54290

$ws.Range("H127").Value = 1362.9546
$ws.Range("I127").Value = 1042.6666
$ws.Range("J127").Value = 2804.25
$ws.Range("K127").Value = 3127.9998
$ws.Range("L127").Value = 8412.75
$ws.Range("M127").Value = 1832.0002
$ws.Range("N127").Value = -18332.75

$ws.Range("H129").Value = 2191.7727
$ws.Range("I129").Value = 1672.1666
$ws.Range("K129").Value = 5016.4998
$ws.Range("M129").Value = -16.4997999999996

$ws.Range("H132").Value = 1167.7073
$ws.Range("I132").Value = 623.0789
$ws.Range("K132").Value = 1869.2367
$ws.Range("M132").Value = 660.7633000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 12989101
$ws.Range("I2").Value = 15152284
$ws.Range("K2").Value = 15152284
$ws.Range("M2").Value = -15152171

$ws.Range("H32").Value = 53204.848
$ws.Range("I32").Value = 60015.668
$ws.Range("K32").Value = 60015.668
$ws.Range("M32").Value = -59728.668

$ws.Range("H35").Value = 1444.1666
$ws.Range("I35").Value = 924.8
$ws.Range("K35").Value = 924.8
$ws.Range("M35").Value = -518.8

$ws.Range("H37").Value = 54500.375
$ws.Range("I37").Value = 0
$ws.Range("J37").Value = 54500.375
$ws.Range("K37").Value = 0
$ws.Range("L37").Value = 54500.375
$ws.Range("N37").Value = -55046.375
$ws.Range("M37").ClearContents()

$ws.Range("H42").Value = 16000
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = 16000
$ws.Range("K42").Value = 0
$ws.Range("L42").Value = 16000
$ws.Range("N42").Value = -16972
$ws.Range("M42").ClearContents()

$ws.Range("H61").Value = 6653.8823
$ws.Range("I61").Value = 5007.0386
$ws.Range("K61").Value = 5007.0386
$ws.Range("M61").Value = -4795.0386

$ws.Range("H63").Value = 8984.08
$ws.Range("I63").Value = 4899
$ws.Range("K63").Value = 4899
$ws.Range("M63").Value = -4213

$ws.Range("H66").Value = 8984.08
$ws.Range("I66").Value = 4899
$ws.Range("K66").Value = 24495
$ws.Range("M66").Value = -21063

$ws.Range("H105").Value = 96708.77
$ws.Range("J105").Value = 96708.77
$ws.Range("L105").Value = 96708.77
$ws.Range("N105").Value = -103696.77

$ws.Range("H116").Value = 12989101
$ws.Range("I116").Value = 15152284
$ws.Range("K116").Value = 15152284
$ws.Range("M116").Value = -15149990

$ws.Range("H122").Value = 3874.75
$ws.Range("I122").Value = 1749.5
$ws.Range("K122").Value = 5248.5
$ws.Range("M122").Value = -2798.5

$ws.Range("H132").Value = 7581.4414
$ws.Range("I132").Value = 6014.85
$ws.Range("J132").Value = 9819.429
$ws.Range("K132").Value = 18044.55
$ws.Range("L132").Value = 29458.287
$ws.Range("M132").Value = -15514.55
$ws.Range("N132").Value = -34518.287

$ws.Range("H136").Value = 6653.8823
$ws.Range("I136").Value = 5007.0386
$ws.Range("K136").Value = 15021.1158
$ws.Range("M136").Value = -12471.1158

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 12989101
$ws.Range("I3").Value = 15152284
$ws.Range("K3").Value = 15152284
$ws.Range("M3").Value = -15152170

$ws.Range("H86").Value = 85066.28999999999
$ws.Range("I86").Value = 1774.55
$ws.Range("K86").Value = 1774.55
$ws.Range("M86").Value = -651.55

$ws.Range("H89").Value = 85066.28999999999
$ws.Range("I89").Value = 1774.55
$ws.Range("K89").Value = 8872.75
$ws.Range("M89").Value = -3256.75

$ws.Range("H134").Value = 6762.4287
$ws.Range("I134").Value = 5234.75
$ws.Range("K134").Value = 15704.25
$ws.Range("M134").Value = -13169.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1020.63635
$ws.Range("J22").Value = 1946
$ws.Range("L22").Value = 1946
$ws.Range("N22").Value = -2646

$ws.Range("H31").Value = 19611488
$ws.Range("I31").Value = 66669610
$ws.Range("J31").Value = 3937.0833
$ws.Range("K31").Value = 66669610
$ws.Range("L31").Value = 3937.0833
$ws.Range("M31").Value = -66669315
$ws.Range("N31").Value = -4527.0833

$ws.Range("H34").Value = 19611488
$ws.Range("I34").Value = 66669610
$ws.Range("J34").Value = 3937.0833
$ws.Range("K34").Value = 66669610
$ws.Range("L34").Value = 3937.0833
$ws.Range("M34").Value = -66669408
$ws.Range("N34").Value = -4341.0833

$ws.Range("H99").Value = 3534.8235
$ws.Range("I99").Value = 3479.875
$ws.Range("K99").Value = 3479.875
$ws.Range("M99").Value = -1981.875

$ws.Range("H107").Value = 549.61536
$ws.Range("I107").Value = 528.75
$ws.Range("J107").Value = 800
$ws.Range("K107").Value = 528.75
$ws.Range("L107").Value = 800
$ws.Range("M107").Value = 1391.25
$ws.Range("N107").Value = -4640

$ws.Range("H122").Value = 4238
$ws.Range("I122").Value = 3578
$ws.Range("K122").Value = 10734
$ws.Range("M122").Value = -8284

$ws.Range("H126").Value = 3534.8235
$ws.Range("I126").Value = 3479.875
$ws.Range("K126").Value = 10439.625
$ws.Range("M126").Value = -7969.625

$ws.Range("H132").Value = 36547.117
$ws.Range("I132").Value = 3612.1538
$ws.Range("K132").Value = 10836.4614
$ws.Range("M132").Value = -8306.4614

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1488797.1
$ws.Range("I4").Value = 895727
$ws.Range("J4").Value = 4651838
$ws.Range("K4").Value = 2687181
$ws.Range("L4").Value = 13955514
$ws.Range("M4").Value = -2687069
$ws.Range("N4").Value = -13955738

$ws.Range("H103").Value = 234.6
$ws.Range("J103").Value = 350.6
$ws.Range("L103").Value = 1051.8
$ws.Range("N103").Value = -2809.8

$ws.Range("H122").Value = 4886.375
$ws.Range("I122").Value = 8051.1113
$ws.Range("J122").Value = 817.4286
$ws.Range("K122").Value = 72460.00169999999
$ws.Range("L122").Value = 7356.8574
$ws.Range("M122").Value = -70010.00169999999
$ws.Range("N122").Value = -12256.8574

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2945.8235
$ws.Range("I80").Value = 2279.5
$ws.Range("J80").Value = 3897.7144
$ws.Range("K80").Value = 2279.5
$ws.Range("L80").Value = 3897.7144
$ws.Range("M80").Value = -1281.5
$ws.Range("N80").Value = -5893.7144

$ws.Range("H83").Value = 2945.8235
$ws.Range("I83").Value = 2279.5
$ws.Range("J83").Value = 3897.7144
$ws.Range("K83").Value = 11397.5
$ws.Range("L83").Value = 19488.572
$ws.Range("M83").Value = -6405.5
$ws.Range("N83").Value = -29472.572

$ws.Range("H113").Value = 7126590
$ws.Range("J113").Value = 19001372
$ws.Range("L113").Value = 19001372
$ws.Range("N113").Value = -19005712

$ws.Range("H114").Value = 70663
$ws.Range("I114").Value = 54000
$ws.Range("J114").Value = 78994.5
$ws.Range("K114").Value = 54000
$ws.Range("L114").Value = 78994.5
$ws.Range("M114").Value = -49661
$ws.Range("N114").Value = -87672.5

$ws.Range("H132").Value = 9434
$ws.Range("I132").Value = 4864.5
$ws.Range("K132").Value = 14593.5
$ws.Range("M132").Value = -12063.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3710.5134
$ws.Range("I22").Value = 2510.375
$ws.Range("K22").Value = 2510.375
$ws.Range("M22").Value = -2215.375

$ws.Range("H27").Value = 3710.5134
$ws.Range("I27").Value = 2510.375
$ws.Range("K27").Value = 2510.375
$ws.Range("M27").Value = -2403.375

$ws.Range("H68").Value = 3172.7273
$ws.Range("I68").Value = 2600
$ws.Range("J68").Value = 3860
$ws.Range("K68").Value = 2600
$ws.Range("L68").Value = 3860
$ws.Range("M68").Value = -1851
$ws.Range("N68").Value = -5358

$ws.Range("H71").Value = 3172.7273
$ws.Range("I71").Value = 2600
$ws.Range("J71").Value = 3860
$ws.Range("K71").Value = 13000
$ws.Range("L71").Value = 19300
$ws.Range("M71").Value = -9256
$ws.Range("N71").Value = -26788

$ws.Range("H132").Value = 4818.4
$ws.Range("I132").Value = 3162.4443
$ws.Range("K132").Value = 9487.332900000001
$ws.Range("M132").Value = -6957.332900000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 5477.4814
$ws.Range("I132").Value = 3957.9473
$ws.Range("K132").Value = 11873.8419
$ws.Range("M132").Value = -9343.841899999999
